# "Reiterable form" rework:
#   - C2 / C3 sample values updated (13 -> 0, 5 -> 50)
#   - the placeholder rows of the form (5:16) are re-heighted to match the
#     data rows above them (18.75 -> 19.5) and the C column placeholders
#     pick up the same bordered / right-aligned numeric look already used
#     by C2:C4, so the whole form reads as one consistent, reiterable block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated sample data -------------------------------------------------
$ws.Range("C2").Value2 = 0
$ws.Range("C3").Value2 = 50

# --- make the reiterable (placeholder) rows match the filled-in rows ----
# Row height: 18.75 -> 19.5, same as rows 1-4.
$ws.Range("A5:A16").EntireRow.RowHeight = 19.5

# Column C placeholders (C5:C16) adopt C4's formatting (bordered, right
# aligned numeric style) instead of their old borderless "general" style.
$ws.Range("C4").Copy()
$ws.Range("C5:C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Output "Form rows re-heighted and C5:C16 restyled to match C2:C4; C2/C3 values updated."
